# Update the crypto price/volume table with freshly scraped values.
# Numeric-looking price strings (e.g. "1.001", "0.08502") are written with
# NumberFormat "@" first so Excel stores them as text instead of coercing
# them to floating point numbers (which would lose trailing zeros / change
# precision), then the style is reset back to "Normal" so no stray
# number-format styling is left on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.093.47'
$ws.Range('D3').Value = '1.918.99'
$ws.Range('E3').Value = '  +2.53%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.63%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.03'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +4.57%  '
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5223'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +2.60%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4084'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.59%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08502'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.60%  '
$ws.Range('E10').Value = '  +2.88%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.128'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.93%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.61'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +10.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.434'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +3.46%  '
$ws.Range('D14').Value = '1.927.58'
$ws.Range('E14').Value = '  +2.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.429'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +2.03%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.001'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.70%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '95.58'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +4.93%  '
$ws.Range('E18').Value = '  +0.85%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06722'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.19%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +3.26%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.000'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.60%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.018'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +1.61%  '
$ws.Range('D23').Value = '30.087.33'
$ws.Range('E23').Value = '  +5.44%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.33'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.222'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +0.51%  '
$ws.Range('D26').Value = '2.147.72'
$ws.Range('E26').Value = '  +2.94%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.53'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.04%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '21.10'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +2.29%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.460'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +1.49%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '129.44'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +1.80%  '
$ws.Range('E31').Value = '  +3.38%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.1056'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +1.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.076'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.640'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.33%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02489'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.22%  '
$ws.Range('E36').Value = '  +0.81%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2215'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.36%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.233'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.197'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.931'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6551'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.88%  '
$ws.Range('B42').Value = 'TrustWalletToken'
$ws.Range('C42').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.249'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.96%  '
$ws.Range('B43').Value = 'Aptos'
$ws.Range('C43').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.65'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +5.06%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.6175'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +2.76%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.29'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.85%  '
$ws.Range('E46').Value = '  +2.23%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.085'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +4.18%  '
$ws.Range('E48').Value = '  +2.73%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '124.64'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +1.92%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.166'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +12.92%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '79.82'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +4.47%  '
